$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.05954295240795
$ws.Range("D2").Value = 1.061020417158375
$ws.Range("E2").Value = 0.992614727750844
$ws.Range("F2").Value = 1.069886259806614
$ws.Range("I2").Value = 1.041176206423329
$ws.Range("J2").Value = 1.064528438569502
$ws.Range("K2").Value = 1.063745612861647
$ws.Range("L2").Value = 0.9955398523335997
$ws.Range("M2").Value = 1.07258757005891
$ws.Range("N2").Value = 1.066040191055389

# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.061446953151594
$ws.Range("D3").Value = 1.062531770625125
$ws.Range("E3").Value = 0.9936372048519299
$ws.Range("F3").Value = 1.071636908074198
$ws.Range("I3").Value = 1.041630489129412
$ws.Range("J3").Value = 1.066081170437255
$ws.Range("K3").Value = 1.065069749205786
$ws.Range("L3").Value = 0.9963617723202687
$ws.Range("M3").Value = 1.074152161457901
$ws.Range("N3").Value = 1.067595127980495

# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.062674897157349
$ws.Range("D4").Value = 1.063505908473114
$ws.Range("E4").Value = 0.9942998659930998
$ws.Range("F4").Value = 1.072766085730645
$ws.Range("I4").Value = 1.041921356050786
$ws.Range("J4").Value = 1.067081583140043
$ws.Range("K4").Value = 1.065922247435578
$ws.Range("L4").Value = 0.9968940712668347
$ws.Range("M4").Value = 1.075160495427255
$ws.Range("N4").Value = 1.06859696138407

# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.063190170824803
$ws.Range("D5").Value = 1.063914539758365
$ws.Range("E5").Value = 0.994578699834602
$ws.Range("F5").Value = 1.073239945860077
$ws.Range("I5").Value = 1.042042903942948
$ws.Range("J5").Value = 1.067501142245324
$ws.Range("K5").Value = 1.066279621363161
$ws.Range("L5").Value = 0.9971179600053012
$ws.Range("M5").Value = 1.075583443374933
$ws.Range("N5").Value = 1.069017116311405

# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.063276632198457
$ws.Range("D6").Value = 1.063983098641492
$ws.Range("E6").Value = 0.994625531979634
$ws.Range("F6").Value = 1.073319459955446
$ws.Range("I6").Value = 1.042063269579498
$ws.Range("J6").Value = 1.06757152907394
$ws.Range("K6").Value = 1.066339566842996
$ws.Range("L6").Value = 0.9971555583673455
$ws.Range("M6").Value = 1.075654402636033
$ws.Range("N6").Value = 1.069087603097391

# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.062681785989616
$ws.Range("D7").Value = 1.063511372124828
$ws.Range("E7").Value = 0.9943035907978918
$ws.Range("F7").Value = 1.072772420773486
$ws.Range("I7").Value = 1.041922983050152
$ws.Range("J7").Value = 1.067087193271971
$ws.Range("K7").Value = 1.065927026657006
$ws.Range("L7").Value = 0.9968970624459044
$ws.Range("M7").Value = 1.075166150610875
$ws.Range("N7").Value = 1.068602579483029

# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.060187274052409
$ws.Range("D8").Value = 1.061531984235503
$ws.Range("E8").Value = 0.9929600610674297
$ws.Range("F8").Value = 1.070478657592704
$ws.Range("I8").Value = 1.041330375566531
$ws.Range("J8").Value = 1.065054094057666
$ws.Range("K8").Value = 1.064194011508639
$ws.Range("L8").Value = 0.9958175282591057
$ws.Range("M8").Value = 1.073117182026451
$ws.Range("N8").Value = 1.06656659303464

# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.055759500535637
$ws.Range("D9").Value = 1.05801417329524
$ws.Range("E9").Value = 0.9906006454969559
$ws.Range("F9").Value = 1.066408309892517
$ws.Range("I9").Value = 1.040262238476178
$ws.Range("J9").Value = 1.061437741452129
$ws.Range("K9").Value = 1.061106571606201
$ws.Range("L9").Value = 0.9939188001724441
$ws.Range("M9").Value = 1.069474781492026
$ws.Range("N9").Value = 1.062945104793601

# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.052784730569532
$ws.Range("D10").Value = 1.055647885812576
$ws.Range("E10").Value = 0.989033133672735
$ws.Range("F10").Value = 1.063674488095137
$ws.Range("I10").Value = 1.039533713669774
$ws.Range("J10").Value = 1.059003045042395
$ws.Range("K10").Value = 1.059024722945813
$ws.Range("L10").Value = 0.9926553831429383
$ws.Range("M10").Value = 1.067024038050622
$ws.Range("N10").Value = 1.060506950835702

# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.051490889294087
$ws.Range("D11").Value = 1.054618034267223
$ws.Range("E11").Value = 0.988355674866747
$ws.Range("F11").Value = 1.062485654266548
$ws.Range("I11").Value = 1.039214272564067
$ws.Range("J11").Value = 1.057942906543361
$ws.Range("K11").Value = 1.058117462679461
$ws.Range("L11").Value = 0.9921088820399291
$ws.Range("M11").Value = 1.06595727039613
$ws.Range("N11").Value = 1.0594453068184

# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.051009410021613
$ws.Range("D12").Value = 1.054234695909902
$ws.Range("E12").Value = 0.9881042295826724
$ws.Range("F12").Value = 1.062043284537584
$ws.Range("I12").Value = 1.039095011866669
$ws.Range("J12").Value = 1.057548216450246
$ws.Range("K12").Value = 1.057779575291489
$ws.Range("L12").Value = 0.9919059725120875
$ws.Range("M12").Value = 1.065560166492821
$ws.Range("N12").Value = 1.059050056220081

# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.051112729606235
$ws.Range("D13").Value = 1.054316960059624
$ws.Range("E13").Value = 0.9881581567098651
$ws.Range("F13").Value = 1.062138210177214
$ws.Range("I13").Value = 1.039120621234977
$ws.Range("J13").Value = 1.057632920222831
$ws.Range("K13").Value = 1.057852093877018
$ws.Range("L13").Value = 0.9919494934313052
$ws.Range("M13").Value = 1.065645385829355
$ws.Range("N13").Value = 1.059134880281739

# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.051451108352849
$ws.Range("D14").Value = 1.054586363946882
$ws.Range("E14").Value = 0.9883348863814464
$ws.Range("F14").Value = 1.062449103967053
$ws.Range("I14").Value = 1.039204426843642
$ws.Range("J14").Value = 1.057910299972669
$ws.Range("K14").Value = 1.05808955107552
$ws.Range("L14").Value = 0.9920921077337197
$ws.Range("M14").Value = 1.065924463305277
$ws.Range("N14").Value = 1.059412653942638

# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.051659476088318
$ws.Range("D15").Value = 1.054752245233959
$ws.Range("E15").Value = 0.9884438009545853
$ws.Range("F15").Value = 1.062640551350529
$ws.Range("I15").Value = 1.039255981716506
$ws.Range("J15").Value = 1.058081081905338
$ws.Range("K15").Value = 1.058235737738825
$ws.Range("L15").Value = 0.9921799884222134
$ws.Range("M15").Value = 1.06609629771429
$ws.Range("N15").Value = 1.05958367840524

# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.052870475300988
$ws.Range("D16").Value = 1.055716121757839
$ws.Range("E16").Value = 0.9890781214508737
$ws.Range("F16").Value = 1.06375327836608
$ws.Range("I16").Value = 1.03955482933946
$ws.Range("J16").Value = 1.059073276856675
$ws.Range("K16").Value = 1.059084810944214
$ws.Range("L16").Value = 0.9926916645766087
$ws.Range("M16").Value = 1.067094716597027
$ws.Range("N16").Value = 1.060577282387214

# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.053628546643669
$ws.Range("D17").Value = 1.05631932130258
$ws.Range("E17").Value = 0.9894763578477731
$ws.Range("F17").Value = 1.064449889511487
$ws.Range("I17").Value = 1.039741216498404
$ws.Range("J17").Value = 1.059694060356259
$ws.Range("K17").Value = 1.059615844895436
$ws.Range("L17").Value = 0.9930127773692701
$ws.Range("M17").Value = 1.06771949010085
$ws.Range("N17").Value = 1.061198947470572

# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.054070164156548
$ws.Range("D18").Value = 1.056670653762091
$ws.Range("E18").Value = 0.9897087662937551
$ws.Range("F18").Value = 1.064855723239279
$ws.Range("I18").Value = 1.039849548975799
$ws.Range("J18").Value = 1.060055585159373
$ws.Range("K18").Value = 1.059925029059184
$ws.Range("L18").Value = 0.9932001317071766
$ws.Range("M18").Value = 1.068083372975701
$ws.Range("N18").Value = 1.061560985680375

# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.054220651345152
$ws.Range("D19").Value = 1.056790364209245
$ws.Range("E19").Value = 0.9897880325774039
$ws.Range("F19").Value = 1.064994019981775
$ws.Range("I19").Value = 1.039886422638499
$ws.Range("J19").Value = 1.06017876011857
$ws.Range("K19").Value = 1.060030358608244
$ws.Range("L19").Value = 0.993264023964098
$ws.Range("M19").Value = 1.068207357178066
$ws.Range("N19").Value = 1.061684335562141

# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.053547270113938
$ws.Range("D20").Value = 1.056254655932402
$ws.Range("E20").Value = 0.9894336180355766
$ws.Range("F20").Value = 1.064375200339275
$ws.Range("I20").Value = 1.039721258691886
$ws.Range("J20").Value = 1.059627515035438
$ws.Range("K20").Value = 1.059558927924895
$ws.Range("L20").Value = 0.9929783193490043
$ws.Range("M20").Value = 1.067652513503507
$ws.Range("N20").Value = 1.061132307647763

# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.051351488995201
$ws.Range("D21").Value = 1.054507053613327
$ws.Range("E21").Value = 0.9882828385668255
$ws.Range("F21").Value = 1.062357575313271
$ws.Range("I21").Value = 1.039179764970033
$ws.Range("J21").Value = 1.057828643771788
$ws.Range("K21").Value = 1.058019650579721
$ws.Range("L21").Value = 0.9920501090198107
$ws.Range("M21").Value = 1.065842305821374
$ws.Range("N21").Value = 1.059330881780585

# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.049965757698136
$ws.Range("D22").Value = 1.05340359569656
$ws.Range("E22").Value = 0.9875604150241496
$ws.Range("F22").Value = 1.061084467753743
$ws.Range("I22").Value = 1.038835796257467
$ws.Range("J22").Value = 1.056692361049828
$ws.Range("K22").Value = 1.057046684269476
$ws.Range("L22").Value = 0.991467000034148
$ws.Range("M22").Value = 1.064699176650749
$ws.Range("N22").Value = 1.058192985406826

# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.050700857535066
$ws.Range("D23").Value = 1.053989008910343
$ws.Range("E23").Value = 0.9879432794636459
$ws.Range("F23").Value = 1.061759804571973
$ws.Range("I23").Value = 1.039018475776534
$ws.Range("J23").Value = 1.057295231857732
$ws.Range("K23").Value = 1.057562967530483
$ws.Range("L23").Value = 0.9917760702887607
$ws.Range("M23").Value = 1.065305650039683
$ws.Range("N23").Value = 1.058796712360428

# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.053583997186266
$ws.Range("D24").Value = 1.056283876985443
$ws.Range("E24").Value = 0.9894529299347241
$ws.Range("F24").Value = 1.064408950667896
$ws.Range("I24").Value = 1.039730277951949
$ws.Range("J24").Value = 1.059657585755161
$ws.Range("K24").Value = 1.05958464798062
$ws.Range("L24").Value = 0.9929938892766438
$ws.Range("M24").Value = 1.067682779003489
$ws.Range("N24").Value = 1.061162421071357

# Row 25
$ws.Range("B25").Value = 1.019999999999999
$ws.Range("C25").Value = 1.056908128114461
$ws.Range("D25").Value = 1.0589272536608
$ws.Range("E25").Value = 0.9912096547607046
$ws.Range("F25").Value = 1.067464078762897
$ws.Range("I25").Value = 1.040541246334985
$ws.Range("J25").Value = 1.062376768582061
$ws.Range("K25").Value = 1.061908831382576
$ws.Range("L25").Value = 0.9944092447426411
$ws.Range("M25").Value = 1.070420312514536
$ws.Range("N25").Value = 1.063885465449766

Write-Output "Applied vm_pu updates for rows 2-25"
